$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Cells.Item(33, 8).Value = 751.8333
$ws.Cells.Item(33, 9).Value = 863.2
$ws.Cells.Item(33, 10).Value = 195
$ws.Cells.Item(33, 11).Value = 863.2
$ws.Cells.Item(33, 12).Value = 195
$ws.Cells.Item(33, 13).Value = -634.2
$ws.Cells.Item(33, 14).Value = -653
# Row 96
$ws.Cells.Item(96, 8).Value = 457.8889
$ws.Cells.Item(96, 9).Value = 457.8889
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 11).Value = 1373.6667
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 13).Value = -0.6666999999999916
# Row 98
$ws.Cells.Item(98, 8).Value = 1815.625
$ws.Cells.Item(98, 9).Value = 1532.1428
$ws.Cells.Item(98, 10).Value = 3800
$ws.Cells.Item(98, 11).Value = 1532.1428
$ws.Cells.Item(98, 12).Value = 3800
$ws.Cells.Item(98, 13).Value = -34.14280000000008
$ws.Cells.Item(98, 14).Value = -6796
# Row 100
$ws.Cells.Item(100, 8).Value = 1476.3636
$ws.Cells.Item(100, 9).Value = 1062.5714
$ws.Cells.Item(100, 10).Value = 2200.5
$ws.Cells.Item(100, 11).Value = 1062.5714
$ws.Cells.Item(100, 12).Value = 2200.5
$ws.Cells.Item(100, 13).Value = -521.5714
$ws.Cells.Item(100, 14).Value = -3282.5
# Row 122
$ws.Cells.Item(122, 8).Value = 1815.625
$ws.Cells.Item(122, 9).Value = 1532.1428
$ws.Cells.Item(122, 10).Value = 3800
$ws.Cells.Item(122, 11).Value = 4596.428400000001
$ws.Cells.Item(122, 12).Value = 11400
$ws.Cells.Item(122, 13).Value = -2146.428400000001
$ws.Cells.Item(122, 14).Value = -16300

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 3779.8
$ws.Cells.Item(45, 9).Value = 1949.5
$ws.Cells.Item(45, 10).Value = 5000
$ws.Cells.Item(45, 11).Value = 1949.5
$ws.Cells.Item(45, 12).Value = 5000
$ws.Cells.Item(45, 13).Value = -1572.5
$ws.Cells.Item(45, 14).Value = -5754
# Row 74
$ws.Cells.Item(74, 8).Value = 68182730
$ws.Cells.Item(74, 9).Value = 107143600
$ws.Cells.Item(74, 10).Value = 1199.5
$ws.Cells.Item(74, 11).Value = 107143600
$ws.Cells.Item(74, 12).Value = 1199.5
$ws.Cells.Item(74, 13).Value = -107142726
$ws.Cells.Item(74, 14).Value = -2947.5
# Row 76
$ws.Cells.Item(76, 8).Value = 199999
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 199999
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 12).Value = 199999
$ws.Cells.Item(76, 14).Value = -200675
# Row 77
$ws.Cells.Item(77, 8).Value = 68182730
$ws.Cells.Item(77, 9).Value = 107143600
$ws.Cells.Item(77, 10).Value = 1199.5
$ws.Cells.Item(77, 11).Value = 535718000
$ws.Cells.Item(77, 12).Value = 5997.5
$ws.Cells.Item(77, 13).Value = -535713632
$ws.Cells.Item(77, 14).Value = -14733.5
# Row 79
$ws.Cells.Item(79, 8).Value = 199999
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 199999
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 12).Value = 199999
$ws.Cells.Item(79, 14).Value = -202339
# Row 102
$ws.Cells.Item(102, 8).Value = 2134.2144
$ws.Cells.Item(102, 9).Value = 2183
$ws.Cells.Item(102, 10).Value = 1500
$ws.Cells.Item(102, 11).Value = 2183
$ws.Cells.Item(102, 12).Value = 1500
$ws.Cells.Item(102, 13).Value = -561
$ws.Cells.Item(102, 14).Value = -4744
# Row 110
$ws.Cells.Item(110, 8).Value = 844.25
$ws.Cells.Item(110, 9).Value = 750.41174
$ws.Cells.Item(110, 10).Value = 1072.1428
$ws.Cells.Item(110, 11).Value = 750.41174
$ws.Cells.Item(110, 12).Value = 1072.1428
$ws.Cells.Item(110, 13).Value = 1294.58826
$ws.Cells.Item(110, 14).Value = -5162.1428
# Row 122
$ws.Cells.Item(122, 8).Value = 6416045
$ws.Cells.Item(122, 9).Value = 9263899
$ws.Cells.Item(122, 10).Value = 8374.75
$ws.Cells.Item(122, 11).Value = 27791697
$ws.Cells.Item(122, 12).Value = 25124.25
$ws.Cells.Item(122, 13).Value = -27789247
$ws.Cells.Item(122, 14).Value = -30024.25

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Cells.Item(22, 8).Value = 1050
$ws.Cells.Item(22, 9).Value = 1000
$ws.Cells.Item(22, 10).Value = 1100
$ws.Cells.Item(22, 11).Value = 1000
$ws.Cells.Item(22, 12).Value = 1100
$ws.Cells.Item(22, 13).Value = -827
$ws.Cells.Item(22, 14).Value = -1446
# Row 99
$ws.Cells.Item(99, 8).Value = 2013.5454
$ws.Cells.Item(99, 9).Value = 1828.375
$ws.Cells.Item(99, 10).Value = 2507.3333
$ws.Cells.Item(99, 11).Value = 1828.375
$ws.Cells.Item(99, 12).Value = 2507.3333
$ws.Cells.Item(99, 13).Value = -330.375
$ws.Cells.Item(99, 14).Value = -5503.3333
# Row 105
$ws.Cells.Item(105, 8).Value = 3758.4285
$ws.Cells.Item(105, 9).Value = 3892.5454
$ws.Cells.Item(105, 10).Value = 3266.6667
$ws.Cells.Item(105, 11).Value = 3892.5454
$ws.Cells.Item(105, 12).Value = 3266.6667
$ws.Cells.Item(105, 13).Value = -2145.5454
$ws.Cells.Item(105, 14).Value = -6760.6667
# Row 107
$ws.Cells.Item(107, 8).Value = 45141028
$ws.Cells.Item(107, 9).Value = 55557100
$ws.Cells.Item(107, 10).Value = 4722.6665
$ws.Cells.Item(107, 11).Value = 55557100
$ws.Cells.Item(107, 12).Value = 4722.6665
$ws.Cells.Item(107, 13).Value = -55555180
$ws.Cells.Item(107, 14).Value = -8562.666499999999
# Row 134
$ws.Cells.Item(134, 8).Value = 2823.2173
$ws.Cells.Item(134, 9).Value = 2522.4443
$ws.Cells.Item(134, 10).Value = 3906
$ws.Cells.Item(134, 11).Value = 7567.3329
$ws.Cells.Item(134, 12).Value = 11718
$ws.Cells.Item(134, 13).Value = -5032.3329
$ws.Cells.Item(134, 14).Value = -16788

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Cells.Item(62, 8).Value = 115000
$ws.Cells.Item(62, 9).Value = 100000
$ws.Cells.Item(62, 10).Value = 120000
$ws.Cells.Item(62, 11).Value = 100000
$ws.Cells.Item(62, 12).Value = 120000
$ws.Cells.Item(62, 13).Value = -99376
$ws.Cells.Item(62, 14).Value = -121248
# Row 65
$ws.Cells.Item(65, 8).Value = 115000
$ws.Cells.Item(65, 9).Value = 100000
$ws.Cells.Item(65, 10).Value = 120000
$ws.Cells.Item(65, 11).Value = 500000
$ws.Cells.Item(65, 12).Value = 600000
$ws.Cells.Item(65, 13).Value = -496880
$ws.Cells.Item(65, 14).Value = -606240
# Row 105
$ws.Cells.Item(105, 8).Value = 1837
$ws.Cells.Item(105, 9).Value = 1699.6666
$ws.Cells.Item(105, 10).Value = 1919.4
$ws.Cells.Item(105, 11).Value = 1699.6666
$ws.Cells.Item(105, 12).Value = 1919.4
$ws.Cells.Item(105, 13).Value = 47.33339999999998
$ws.Cells.Item(105, 14).Value = -5413.4
# Row 141
$ws.Cells.Item(141, 8).Value = 102318.11
$ws.Cells.Item(141, 9).Value = 90765.664
$ws.Cells.Item(141, 10).Value = 103011.26
$ws.Cells.Item(141, 11).Value = 90765.664
$ws.Cells.Item(141, 12).Value = 103011.26
$ws.Cells.Item(141, 13).Value = -85585.664
$ws.Cells.Item(141, 14).Value = -113371.26

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Cells.Item(8, 8).Value = 210.66667
$ws.Cells.Item(8, 9).Value = 210.66667
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 632.00001
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = -493.00001
# Row 68
$ws.Cells.Item(68, 8).Value = 1406.25
$ws.Cells.Item(68, 9).Value = 825
$ws.Cells.Item(68, 10).Value = 1600
$ws.Cells.Item(68, 11).Value = 2475
$ws.Cells.Item(68, 12).Value = 4800
$ws.Cells.Item(68, 13).Value = -1664
$ws.Cells.Item(68, 14).Value = -6422
# Row 71
$ws.Cells.Item(71, 8).Value = 1406.25
$ws.Cells.Item(71, 9).Value = 825
$ws.Cells.Item(71, 10).Value = 1600
$ws.Cells.Item(71, 11).Value = 7425
$ws.Cells.Item(71, 12).Value = 14400
$ws.Cells.Item(71, 13).Value = -3369
$ws.Cells.Item(71, 14).Value = -22512
# Row 103
$ws.Cells.Item(103, 8).Value = 469.42856
$ws.Cells.Item(103, 9).Value = 445.5
$ws.Cells.Item(103, 10).Value = 479
$ws.Cells.Item(103, 11).Value = 1336.5
$ws.Cells.Item(103, 12).Value = 1437
$ws.Cells.Item(103, 13).Value = -457.5
$ws.Cells.Item(103, 14).Value = -3195
# Row 141
$ws.Cells.Item(141, 8).Value = 3355.9333
$ws.Cells.Item(141, 9).Value = 2333.8462
$ws.Cells.Item(141, 10).Value = 9999.5
$ws.Cells.Item(141, 11).Value = 7001.5386
$ws.Cells.Item(141, 12).Value = 29998.5
$ws.Cells.Item(141, 13).Value = -1821.5386
$ws.Cells.Item(141, 14).Value = -40358.5

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 128093.5
$ws.Cells.Item(80, 9).Value = 253812.25
$ws.Cells.Item(80, 10).Value = 2374.75
$ws.Cells.Item(80, 11).Value = 253812.25
$ws.Cells.Item(80, 12).Value = 2374.75
$ws.Cells.Item(80, 13).Value = -252814.25
$ws.Cells.Item(80, 14).Value = -4370.75
# Row 83
$ws.Cells.Item(83, 8).Value = 128093.5
$ws.Cells.Item(83, 9).Value = 253812.25
$ws.Cells.Item(83, 10).Value = 2374.75
$ws.Cells.Item(83, 11).Value = 1269061.25
$ws.Cells.Item(83, 12).Value = 11873.75
$ws.Cells.Item(83, 13).Value = -1264069.25
$ws.Cells.Item(83, 14).Value = -21857.75
# Row 102
$ws.Cells.Item(102, 8).Value = 29422936
$ws.Cells.Item(102, 9).Value = 50012996
$ws.Cells.Item(102, 10).Value = 8567
$ws.Cells.Item(102, 11).Value = 50012996
$ws.Cells.Item(102, 12).Value = 8567
$ws.Cells.Item(102, 13).Value = -50011374
$ws.Cells.Item(102, 14).Value = -11811
# Row 126
$ws.Cells.Item(126, 8).Value = 4855.5
$ws.Cells.Item(126, 9).Value = 2065.8333
$ws.Cells.Item(126, 10).Value = 6529.3
$ws.Cells.Item(126, 11).Value = 6197.499899999999
$ws.Cells.Item(126, 12).Value = 19587.9
$ws.Cells.Item(126, 13).Value = -3727.499899999999
$ws.Cells.Item(126, 14).Value = -24527.9

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Cells.Item(68, 8).Value = 3724.9333
$ws.Cells.Item(68, 9).Value = 2845.9048
$ws.Cells.Item(68, 10).Value = 5776
$ws.Cells.Item(68, 11).Value = 2845.9048
$ws.Cells.Item(68, 12).Value = 5776
$ws.Cells.Item(68, 13).Value = -2096.9048
$ws.Cells.Item(68, 14).Value = -7274
# Row 71
$ws.Cells.Item(71, 8).Value = 3724.9333
$ws.Cells.Item(71, 9).Value = 2845.9048
$ws.Cells.Item(71, 10).Value = 5776
$ws.Cells.Item(71, 11).Value = 14229.524
$ws.Cells.Item(71, 12).Value = 28880
$ws.Cells.Item(71, 13).Value = -10485.524
$ws.Cells.Item(71, 14).Value = -36368
# Row 93
$ws.Cells.Item(93, 8).Value = 1160.6154
$ws.Cells.Item(93, 9).Value = 1003.7895
$ws.Cells.Item(93, 10).Value = 1586.2858
$ws.Cells.Item(93, 11).Value = 1003.7895
$ws.Cells.Item(93, 12).Value = 1586.2858
$ws.Cells.Item(93, 13).Value = 244.2105
$ws.Cells.Item(93, 14).Value = -4082.2858

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 5021.8237
$ws.Cells.Item(81, 9).Value = 4523.25
$ws.Cells.Item(81, 10).Value = 5465
$ws.Cells.Item(81, 11).Value = 9046.5
$ws.Cells.Item(81, 12).Value = 10930
$ws.Cells.Item(81, 13).Value = -7985.5
$ws.Cells.Item(81, 14).Value = -13052
# Row 84
$ws.Cells.Item(84, 8).Value = 5021.8237
$ws.Cells.Item(84, 9).Value = 4523.25
$ws.Cells.Item(84, 10).Value = 5465
$ws.Cells.Item(84, 11).Value = 45232.5
$ws.Cells.Item(84, 12).Value = 54650
$ws.Cells.Item(84, 13).Value = -39928.5
$ws.Cells.Item(84, 14).Value = -65258
# Row 96
$ws.Cells.Item(96, 8).Value = 6659.6
$ws.Cells.Item(96, 9).Value = 1999
$ws.Cells.Item(96, 10).Value = 7824.75
$ws.Cells.Item(96, 11).Value = 1999
$ws.Cells.Item(96, 12).Value = 7824.75
$ws.Cells.Item(96, 13).Value = -626
$ws.Cells.Item(96, 14).Value = -10570.75
# Row 126
$ws.Cells.Item(126, 8).Value = 1652.3334
$ws.Cells.Item(126, 9).Value = 1584.8
$ws.Cells.Item(126, 10).Value = 1990
$ws.Cells.Item(126, 11).Value = 4754.4
$ws.Cells.Item(126, 12).Value = 5970
$ws.Cells.Item(126, 13).Value = -2284.4
$ws.Cells.Item(126, 14).Value = -10910

Write-Host "Applied all changes"